$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.040.46"
$ws.Range("E2").Value = "  +2.77%  "

# Row 3
$ws.Range("D3").Value = "2.610.20"
$ws.Range("E3").Value = "  +1.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.89"
$ws.Range("E5").Value = "  +4.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.45"
$ws.Range("E6").Value = "  +0.80%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  +0.55%  "

# Row 9
$ws.Range("D9").Value = "2.636.18"
$ws.Range("E9").Value = "  +1.69%  "

# Row 10
$ws.Range("E10").Value = "  -2.36%  "

# Row 11
$ws.Range("E11").Value = "  +1.63%  "

# Row 12
$ws.Range("E12").Value = "  -6.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.366"
$ws.Range("E13").Value = "  +3.15%  "

# Row 14
$ws.Range("D14").Value = "3.072.61"
$ws.Range("E14").Value = "  +1.13%  "

# Row 15
$ws.Range("D15").Value = "60.943.54"
$ws.Range("E15").Value = "  +2.62%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.27"
$ws.Range("E16").Value = "  +0.80%  "

# Row 17
$ws.Range("E17").Value = "  +3.58%  "

# Row 18
$ws.Range("D18").Value = "2.622.92"
$ws.Range("E18").Value = "  +1.44%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.33"
$ws.Range("E19").Value = "  +9.25%  "

# Row 20
$ws.Range("E20").Value = "  +1.85%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.95"
$ws.Range("E21").Value = "  +3.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.92"
$ws.Range("E22").Value = "  +7.13%  "

# Row 23
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.526"
$ws.Range("E24").Value = "  +11.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.55"
$ws.Range("E25").Value = "  +0.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("E27").Value = "  -0.02%  "

# Row 28
$ws.Range("E28").Value = "  +4.31%  "

# Row 29
$ws.Range("E29").Value = "  +2.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.87"
$ws.Range("E30").Value = "  +12.19%  "

# Row 31
$ws.Range("E31").Value = "  +3.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.45"
$ws.Range("E33").Value = "  +1.88%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.55"
$ws.Range("E34").Value = "  +2.34%  "

# Row 35
$ws.Range("E35").Value = "  +4.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.978"
$ws.Range("E36").Value = "  +8.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").Value = "  +6.45%  "

# Row 38
$ws.Range("E38").Value = "  +8.09%  "

# Row 39
$ws.Range("E39").Value = "  +1.40%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.87"
$ws.Range("E40").Value = "  +5.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.845"
$ws.Range("E41").Value = "  -1.90%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "298.27"
$ws.Range("E42").Value = "  +1.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.56"
$ws.Range("E43").Value = "  +0.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0987"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.91"
$ws.Range("E46").Value = "  +4.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.606"
$ws.Range("E47").Value = "  +2.12%  "

# Row 48
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0548"
$ws.Range("E48").Value = "  +2.73%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.97"
$ws.Range("E49").Value = "  +9.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0241"
$ws.Range("E50").Value = "  +2.95%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.87"
$ws.Range("E51").Value = "  +5.60%  "
